# Commit: Sun, May 24, 2020 1:04:44 PM
#
# This edit does two things:
#
# 1. Re-points the three tables (on the slides that hold the "Table_0"
#    style) to a different table style GUID.
#
# 2. Swaps the presentation's applied theme colours: the deck currently
#    renders with the "Integral / Red Violet" colour scheme; after the
#    edit it renders with the standard "Office" colour scheme (the
#    scheme that used to live unused in the side theme part). We apply
#    this by writing the 12 "Office" theme colours onto the live/used
#    theme through the Slide.ThemeColorScheme object, which is the
#    supported COM surface for editing theme colours in place (keeps
#    fonts/format scheme, only swaps the 12-colour palette).

$p = $ppt.ActivePresentation

# --- 1. Table style swap -------------------------------------------------
$newTableStyle = "{3B447038-8047-4593-ABF2-F7A0CC656A79}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colour swap (Integral/Red Violet -> Office) ---------------
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le $tcs.Count; $k++) {
    $tcs.Item($k).RGB = $officeColors[$k - 1]
}
